$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "missing" count for row 26 (C26 = 37) from 2 to 4.
$ws.Range("I26").Value = 4

# Restore the green highlight fill on F26:G26 (matches the rest of the
# F:G column body, fill color FF99FF66) which had been cleared to no fill.
$ws.Range("F26:G26").Interior.Color = 6750105

# Move the active selection to A3.
$ws.Range("A3").Select()
